$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2122
$ws.Range("E2").Value = 31
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = 42
$ws.Range("H2").Value = 25
$ws.Range("I2").Value = 21
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 3425
$ws.Range("L2").Value = 1701
$ws.Range("M2").Value = 1725
$ws.Range("N2").Value = 1558
$ws.Range("O2").Value = 167
$ws.Range("P2").Value = 392
$ws.Range("Q2").Value = 91
$ws.Range("R2").Value = -292
$ws.Range("S2").Value = 197
$ws.Range("T2").Value = 284
$ws.Range("U2").Value = -193
$ws.Range("V2").Value = 789
$ws.Range("W2").Value = 1.44
$ws.Range("X2").Value = 1.19
$ws.Range("Y2").Value = 1.34
$ws.Range("Z2").Value = 0.77
$ws.Range("AA2").Value = 98.59
$ws.Range("AB2").Value = 552.15
$ws.Range("AC2").Value = 26
$ws.Range("AD2").Value = 67.23
$ws.Range("AE2").Value = 3434
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 2.81
$ws.Range("AI2").Value = 109.28
$ws.Range("AJ2").Value = 78389202

# Row 3
$ws.Range("D3").Value = 2084
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 29
$ws.Range("G3").Value = 59
$ws.Range("H3").Value = 42
$ws.Range("I3").Value = 37
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 3028
$ws.Range("L3").Value = 1286
$ws.Range("M3").Value = 1742
$ws.Range("N3").Value = 1572
$ws.Range("O3").Value = 170
$ws.Range("P3").Value = 392
$ws.Range("Q3").Value = 51
$ws.Range("R3").Value = 14
$ws.Range("S3").Value = -67
$ws.Range("T3").Value = 46
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 749
$ws.Range("W3").Value = 1.41
$ws.Range("X3").Value = 1.99
$ws.Range("Y3").Value = 2.35
$ws.Range("Z3").Value = 1.29
$ws.Range("AA3").Value = 73.87
$ws.Range("AB3").Value = 554.47
$ws.Range("AC3").Value = 47
$ws.Range("AD3").Value = 49.35
$ws.Range("AE3").Value = 3464
$ws.Range("AF3").Value = 0.67
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 2.16
$ws.Range("AI3").Value = 61.68
$ws.Range("AJ3").Value = 78389202

# Row 4
$ws.Range("D4").Value = 1859
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3008
$ws.Range("L4").Value = 1280
$ws.Range("M4").Value = 1728
$ws.Range("N4").Value = 1727
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 429
$ws.Range("Q4").Value = 94
$ws.Range("R4").Value = -69
$ws.Range("S4").Value = 18
$ws.Range("T4").Value = 94
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 798
$ws.Range("W4").Value = 0.28
$ws.Range("X4").Value = 0.18
$ws.Range("Y4").Value = 0.16
$ws.Range("Z4").Value = 0.11
$ws.Range("AA4").Value = 74.11
$ws.Range("AB4").Value = 548.26
$ws.Range("AC4").Value = 3
$ws.Range("AD4").Value = 493.17
$ws.Range("AE4").Value = 3426
$ws.Range("AF4").Value = 0.46
$ws.Range("AG4").Value = 30
$ws.Range("AH4").Value = 1.9
$ws.Range("AI4").Value = 564.35
$ws.Range("AJ4").Value = 85826509

# Row 5
$ws.Range("D5").Value = 1774
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = -56
$ws.Range("H5").Value = -42
$ws.Range("I5").Value = -42
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2981
$ws.Range("L5").Value = 1311
$ws.Range("M5").Value = 1670
$ws.Range("N5").Value = 1669
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 429
$ws.Range("Q5").Value = -20
$ws.Range("R5").Value = -62
$ws.Range("S5").Value = 20
$ws.Range("T5").Value = 53
$ws.Range("U5").Value = -73
$ws.Range("V5").Value = 830
$ws.Range("W5").Value = 0.61
$ws.Range("X5").Value = -2.36
$ws.Range("Y5").Value = -2.48
$ws.Range("Z5").Value = -1.4
$ws.Range("AA5").Value = 78.48
$ws.Range("AB5").Value = 535.49
$ws.Range("AC5").Value = -49
$ws.Range("AD5").Value = -25.32
$ws.Range("AE5").Value = 3312
$ws.Range("AF5").Value = 0.37
$ws.Range("AG5").Value = 10
$ws.Range("AH5").Value = 0.81
$ws.Range("AI5").Value = -11.99
$ws.Range("AJ5").Value = 85826509

# Row 6
$ws.Range("D6").Value = 1965
$ws.Range("E6").Value = -17
$ws.Range("F6").Value = -17
$ws.Range("G6").Value = -59
$ws.Range("H6").Value = -52
$ws.Range("I6").Value = -52
$ws.Range("K6").Value = 2930
$ws.Range("L6").Value = 1266
$ws.Range("M6").Value = 1664
$ws.Range("N6").Value = 1663
$ws.Range("P6").Value = 429
$ws.Range("Q6").Value = 79
$ws.Range("R6").Value = -59
$ws.Range("S6").Value = 42
$ws.Range("T6").Value = 61
$ws.Range("U6").Value = 18
$ws.Range("V6").Value = 820
$ws.Range("W6").Value = -0.87
$ws.Range("X6").Value = -2.65
$ws.Range("Y6").Value = -3.13
$ws.Range("Z6").Value = -1.76
$ws.Range("AA6").Value = 76.07
$ws.Range("AB6").Value = 522.31
$ws.Range("AC6").Value = -61
$ws.Range("AD6").Value = -17.62
$ws.Range("AE6").Value = 3078
$ws.Range("AF6").Value = 0.35
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 0.93
$ws.Range("AI6").Value = -10.36
$ws.Range("AJ6").Value = 85826509

# Clear all data cells (except A/B/C) for rows 7-9 per the diff
$ws.Range("D7:AJ9").ClearContents()

Write-Host "Edit applied"